$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------------
# Sheet 1: LP1912
# ---------------------------------------------------------------------------
$ws1 = $wb.Worksheets.Item("LP1912")

$ws1.Range("A2").Value = "Última actualización: 01:59:40"
$ws1.Range("A3").Value = "Total filas: 6"

# Insert a new row before the existing row 8 (old row 8 shifts down to row 9)
$ws1.Rows.Item(8).Insert()
$ws1.Range("A8").Value = "01:59:40"
$ws1.Range("B8").Value = "02:00"
$ws1.Range("C8").Value = "14_ABASTO"
$ws1.Range("D8").Value = 1
$ws1.Range("E8").Value = "LP1912"

# Append two new rows at the end (rows 10 and 11)
$ws1.Range("A10").Value = "01:59:40"
$ws1.Range("B10").Value = "03:06"
$ws1.Range("C10").Value = "215_ALUAR"
$ws1.Range("D10").Value = 67
$ws1.Range("E10").Value = "LP1912"

$ws1.Range("A11").Value = "01:59:40"
$ws1.Range("B11").Value = "03:50"
$ws1.Range("C11").Value = "14_ABASTO"
$ws1.Range("D11").Value = 111
$ws1.Range("E11").Value = "LP1912"

# ---------------------------------------------------------------------------
# Sheet 2: LP1912-215
# ---------------------------------------------------------------------------
$ws2 = $wb.Worksheets.Item("LP1912-215")

$ws2.Range("A2").Value = "Última actualización: 01:59:40"
$ws2.Range("A3").Value = "Total filas: 3"

# Append a new row at the end (row 8)
$ws2.Range("A8").Value = "01:59:40"
$ws2.Range("B8").Value = "03:06"
$ws2.Range("C8").Value = "215_ALUAR"
$ws2.Range("D8").Value = 67
$ws2.Range("E8").Value = "LP1912"

# ---------------------------------------------------------------------------
# Sheet 3: 6203-6173
# ---------------------------------------------------------------------------
$ws3 = $wb.Worksheets.Item("6203-6173")

$ws3.Range("A2").Value = "Última actualización: 01:59:40"
